# Adding data buoy node bill of materials
# Insert a new "AA Battery" line (qty 3) above the existing "Breadboard" row,
# pushing Breadboard / Cables and Jumpers / PLA filament down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 15 (shifts 15-17 down to 16-18)
$ws.Rows(15).Insert() | Out-Null

# Populate the newly inserted row
$ws.Range("B15").Value = "AA Battery"
$ws.Range("C15").Value = 3

# Match the cursor position left behind by the edit
$ws.Range("B16").Select() | Out-Null
